# Auto-committed style maintenance edit on the DbLayouts "NegFinAcct" sheet.
#
# The DB-layout documentation sheet "DBD" lists one row per column of the
# NegFinAcct table. Two of those rows ("CreateDate" / LastUpdate") had their
# 型態 (data type) documented as "DATE"; the table was actually changed to use
# a TIMESTAMP column, so the layout sheet is corrected to say "TIMESTAMP"
# instead of "DATE" for those two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# Row 18 -> CreateDate  (建檔日期時間)
# Row 20 -> LastUpdate  (最後更新日期時間)
# Column D holds the 形態 (data type) for each field.
$ws.Range("D18").Value = "TIMESTAMP"
$ws.Range("D20").Value = "TIMESTAMP"

# Leave the selection where the edit ended, matching where the author's
# Excel session left the cursor after making this change.
$ws.Range("D20").Select()
